$d = $word.ActiveDocument

# --- Merge split math runs / remove stray bookmark (paragraph content replacements) ---
$p = $d.Paragraphs.Item(2)
$null = $p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="2A516141" w14:textId="5C97CEF2" w:rsidR="007B5560" w:rsidRPr="00DF5314" w:rsidRDefault="00E80C3E" w:rsidP="00D92330"><w:pPr><w:spacing w:after="0" w:line="480" w:lineRule="auto"/><w:contextualSpacing/><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/><w:i/></w:rPr></w:pPr><m:oMathPara><m:oMathParaPr><m:jc m:val="left"/></m:oMathParaPr><m:oMath><m:sSub><m:sSubPr><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:i/></w:rPr></m:ctrlPr></m:sSubPr><m:e><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>x</m:t></m:r></m:e><m:sub><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>1</m:t></m:r></m:sub></m:sSub><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>=size;</m:t></m:r><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/></w:rPr><m:t xml:space="preserve"> </m:t></m:r></m:oMath></m:oMathPara></w:p>')

$p = $d.Paragraphs.Item(4)
$null = $p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="2487C7D3" w14:textId="46EAF51A" w:rsidR="00DF5314" w:rsidRPr="00DF5314" w:rsidRDefault="00E80C3E" w:rsidP="00D92330"><w:pPr><w:spacing w:after="0" w:line="480" w:lineRule="auto"/><w:contextualSpacing/><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/><w:i/></w:rPr></w:pPr><m:oMathPara><m:oMathParaPr><m:jc m:val="left"/></m:oMathParaPr><m:oMath><m:sSub><m:sSubPr><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:i/></w:rPr></m:ctrlPr></m:sSubPr><m:e><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>x</m:t></m:r></m:e><m:sub><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>2</m:t></m:r></m:sub></m:sSub><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>=3*size</m:t></m:r></m:oMath></m:oMathPara></w:p>')

$p = $d.Paragraphs.Item(6)
$null = $p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="3AAE3DF1" w14:textId="3EC9C275" w:rsidR="00DF5314" w:rsidRPr="00DF5314" w:rsidRDefault="00E80C3E" w:rsidP="00D92330"><w:pPr><w:spacing w:after="0" w:line="480" w:lineRule="auto"/><w:contextualSpacing/><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/><w:i/></w:rPr></w:pPr><m:oMathPara><m:oMathParaPr><m:jc m:val="left"/></m:oMathParaPr><m:oMath><m:sSub><m:sSubPr><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:i/></w:rPr></m:ctrlPr></m:sSubPr><m:e><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>x</m:t></m:r></m:e><m:sub><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>3</m:t></m:r></m:sub></m:sSub><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>=4*size;</m:t></m:r></m:oMath></m:oMathPara></w:p>')

$p = $d.Paragraphs.Item(7)
$null = $p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="498B0F02" w14:textId="5D35DE71" w:rsidR="00DF5314" w:rsidRPr="00DF5314" w:rsidRDefault="00E80C3E" w:rsidP="00D92330"><w:pPr><w:spacing w:after="0" w:line="480" w:lineRule="auto"/><w:contextualSpacing/><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/><w:i/></w:rPr></w:pPr><m:oMathPara><m:oMathParaPr><m:jc m:val="left"/></m:oMathParaPr><m:oMath><m:sSub><m:sSubPr><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:i/></w:rPr></m:ctrlPr></m:sSubPr><m:e><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>y</m:t></m:r></m:e><m:sub><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>3</m:t></m:r></m:sub></m:sSub><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>=1.73206*size;</m:t></m:r></m:oMath></m:oMathPara></w:p>')

$p = $d.Paragraphs.Item(8)
$null = $p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="337F1945" w14:textId="76A36A04" w:rsidR="00DF5314" w:rsidRPr="00DF5314" w:rsidRDefault="00E80C3E" w:rsidP="00D92330"><w:pPr><w:spacing w:after="0" w:line="480" w:lineRule="auto"/><w:contextualSpacing/><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/><w:i/></w:rPr></w:pPr><m:oMathPara><m:oMathParaPr><m:jc m:val="left"/></m:oMathParaPr><m:oMath><m:sSub><m:sSubPr><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:i/></w:rPr></m:ctrlPr></m:sSubPr><m:e><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>x</m:t></m:r></m:e><m:sub><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>4</m:t></m:r></m:sub></m:sSub><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>=3*size;</m:t></m:r></m:oMath></m:oMathPara></w:p>')

$p = $d.Paragraphs.Item(9)
$null = $p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="2DF93725" w14:textId="76271AE7" w:rsidR="00DF5314" w:rsidRPr="00DF5314" w:rsidRDefault="00E80C3E" w:rsidP="00D92330"><w:pPr><w:spacing w:after="0" w:line="480" w:lineRule="auto"/><w:contextualSpacing/><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/><w:i/></w:rPr></w:pPr><m:oMathPara><m:oMathParaPr><m:jc m:val="left"/></m:oMathParaPr><m:oMath><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/></w:rPr><m:t xml:space="preserve"> </m:t></m:r><m:sSub><m:sSubPr><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:i/></w:rPr></m:ctrlPr></m:sSubPr><m:e><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>y</m:t></m:r></m:e><m:sub><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>4</m:t></m:r></m:sub></m:sSub><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>=3.46412*size;</m:t></m:r></m:oMath></m:oMathPara></w:p>')

$p = $d.Paragraphs.Item(10)
$null = $p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="2236E233" w14:textId="55400F12" w:rsidR="00DF5314" w:rsidRPr="00DF5314" w:rsidRDefault="00E80C3E" w:rsidP="00D92330"><w:pPr><w:spacing w:after="0" w:line="480" w:lineRule="auto"/><w:contextualSpacing/><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/><w:i/></w:rPr></w:pPr><m:oMathPara><m:oMathParaPr><m:jc m:val="left"/></m:oMathParaPr><m:oMath><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/></w:rPr><m:t xml:space="preserve"> </m:t></m:r><m:sSub><m:sSubPr><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:i/></w:rPr></m:ctrlPr></m:sSubPr><m:e><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>x</m:t></m:r></m:e><m:sub><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>5</m:t></m:r></m:sub></m:sSub><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>=size;</m:t></m:r></m:oMath></m:oMathPara></w:p>')

$p = $d.Paragraphs.Item(11)
$null = $p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="2E9082CA" w14:textId="158F08D2" w:rsidR="00DF5314" w:rsidRPr="00DF5314" w:rsidRDefault="00E80C3E" w:rsidP="00D92330"><w:pPr><w:spacing w:after="0" w:line="480" w:lineRule="auto"/><w:contextualSpacing/><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/><w:i/></w:rPr></w:pPr><m:oMathPara><m:oMathParaPr><m:jc m:val="left"/></m:oMathParaPr><m:oMath><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/></w:rPr><m:t xml:space="preserve"> </m:t></m:r><m:sSub><m:sSubPr><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:i/></w:rPr></m:ctrlPr></m:sSubPr><m:e><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>y</m:t></m:r></m:e><m:sub><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>5</m:t></m:r></m:sub></m:sSub><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>=3.46412*size;</m:t></m:r></m:oMath></m:oMathPara></w:p>')

$p = $d.Paragraphs.Item(13)
$null = $p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="668F0392" w14:textId="2B522D9E" w:rsidR="00DF5314" w:rsidRPr="00AA3BD9" w:rsidRDefault="00E80C3E" w:rsidP="00D92330"><w:pPr><w:spacing w:after="0" w:line="480" w:lineRule="auto"/><w:contextualSpacing/><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/><w:i/></w:rPr></w:pPr><m:oMathPara><m:oMathParaPr><m:jc m:val="left"/></m:oMathParaPr><m:oMath><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/></w:rPr><m:t xml:space="preserve"> </m:t></m:r><m:sSub><m:sSubPr><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:i/></w:rPr></m:ctrlPr></m:sSubPr><m:e><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>y</m:t></m:r></m:e><m:sub><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>6</m:t></m:r></m:sub></m:sSub><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>=1.73206*size;</m:t></m:r></m:oMath></m:oMathPara></w:p>')

$p = $d.Paragraphs.Item(15)
$null = $p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="6D41C6BB" w14:textId="256BCCC3" w:rsidR="009560AC" w:rsidRPr="009560AC" w:rsidRDefault="00DF5314" w:rsidP="00D92330"><w:pPr><w:spacing w:after="0" w:line="480" w:lineRule="auto"/><w:contextualSpacing/><w:rPr><w:b/><w:bCs/><w:iCs/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r w:rsidRPr="009560AC"><w:rPr><w:b/><w:bCs/><w:iCs/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Hexagon </w:t></w:r><w:r w:rsidR="00AB3490"><w:rPr><w:b/><w:bCs/><w:iCs/><w:lang w:val="en-US"/></w:rPr><w:t>Pivots</w:t></w:r></w:p>')

$p = $d.Paragraphs.Item(16)
$null = $p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="5FDD47D9" w14:textId="17D9457B" w:rsidR="001A0C3E" w:rsidRPr="00DF5314" w:rsidRDefault="001A0C3E" w:rsidP="001A0C3E"><w:pPr><w:spacing w:after="0" w:line="480" w:lineRule="auto"/><w:contextualSpacing/><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/><w:i/></w:rPr></w:pPr><m:oMathPara><m:oMathParaPr><m:jc m:val="left"/></m:oMathParaPr><m:oMath><m:sSub><m:sSubPr><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:i/></w:rPr></m:ctrlPr></m:sSubPr><m:e><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>Offset</m:t></m:r></m:e><m:sub><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>X</m:t></m:r></m:sub></m:sSub><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>=</m:t></m:r><m:d><m:dPr><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:i/></w:rPr></m:ctrlPr></m:dPr><m:e><m:sSub><m:sSubPr><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:i/></w:rPr></m:ctrlPr></m:sSubPr><m:e><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:lang w:val="en-US"/></w:rPr><m:t>ElementIndex</m:t></m:r></m:e><m:sub><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>X</m:t></m:r></m:sub></m:sSub><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>-1</m:t></m:r></m:e></m:d><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>*(3*size+0.866025*interval);</m:t></m:r></m:oMath></m:oMathPara></w:p>')

$p = $d.Paragraphs.Item(17)
$null = $p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="2B0A0B86" w14:textId="3E2B228D" w:rsidR="007766B9" w:rsidRPr="00DF5314" w:rsidRDefault="00FF24D2" w:rsidP="00D92330"><w:pPr><w:spacing w:after="0" w:line="480" w:lineRule="auto"/><w:contextualSpacing/><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/><w:i/></w:rPr></w:pPr><m:oMathPara><m:oMathParaPr><m:jc m:val="left"/></m:oMathParaPr><m:oMath><m:sSub><m:sSubPr><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:i/></w:rPr></m:ctrlPr></m:sSubPr><m:e><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>Offset</m:t></m:r></m:e><m:sub><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>Y</m:t></m:r></m:sub></m:sSub><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>=</m:t></m:r><m:d><m:dPr><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:i/></w:rPr></m:ctrlPr></m:dPr><m:e><m:sSub><m:sSubPr><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:i/></w:rPr></m:ctrlPr></m:sSubPr><m:e><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:lang w:val="en-US"/></w:rPr><m:t>ElementIndex</m:t></m:r></m:e><m:sub><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>Y</m:t></m:r></m:sub></m:sSub><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>-1</m:t></m:r></m:e></m:d><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>*(1.73206*size+0.5*interval);</m:t></m:r></m:oMath></m:oMathPara></w:p>')

# --- Append new "Window Size" paragraphs before the final trailing empty paragraph ---
$last = $d.Paragraphs.Last
$null = $last.Range.InsertParagraphBefore()
$null = $last.Range.InsertParagraphBefore()
$null = $last.Range.InsertParagraphBefore()
$null = $last.Range.InsertParagraphBefore()
$null = $last.Range.InsertParagraphBefore()

$n = $d.Paragraphs.Count
# The 5 new placeholder paragraphs are now immediately before the last (trailing empty) paragraph
$target = $d.Paragraphs.Item($n - 5)
$null = $target.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math"><w:pPr><w:spacing w:after="0" w:line="480" w:lineRule="auto"/><w:contextualSpacing/><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/><w:i/></w:rPr></w:pPr></w:p>')

$target = $d.Paragraphs.Item($n - 4)
$null = $target.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math"><w:pPr><w:spacing w:after="0" w:line="480" w:lineRule="auto"/><w:contextualSpacing/><w:rPr><w:b/><w:bCs/><w:iCs/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:iCs/><w:lang w:val="en-US"/></w:rPr><w:t>Window Size</w:t></w:r></w:p>')

$target = $d.Paragraphs.Item($n - 3)
$null = $target.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math"><w:pPr><w:spacing w:after="0" w:line="480" w:lineRule="auto"/><w:contextualSpacing/><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/><w:i/></w:rPr></w:pPr><m:oMathPara><m:oMathParaPr><m:jc m:val="left"/></m:oMathParaPr><m:oMath><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>Width</m:t></m:r><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>=</m:t></m:r><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>10</m:t></m:r><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>*size+</m:t></m:r><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>1.73205</m:t></m:r><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>*interval;</m:t></m:r></m:oMath></m:oMathPara><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>')

$target = $d.Paragraphs.Item($n - 2)
$null = $target.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math"><w:pPr><w:spacing w:after="0" w:line="480" w:lineRule="auto"/><w:contextualSpacing/><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/><w:i/></w:rPr></w:pPr><m:oMathPara><m:oMathParaPr><m:jc m:val="left"/></m:oMathParaPr><m:oMath><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>Height</m:t></m:r><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>=</m:t></m:r><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>10.39236</m:t></m:r><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>*size+</m:t></m:r><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>2</m:t></m:r><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>*interval</m:t></m:r><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>;</m:t></m:r></m:oMath></m:oMathPara></w:p>')

$target = $d.Paragraphs.Item($n - 1)
$null = $target.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math"><w:pPr><w:spacing w:after="0" w:line="480" w:lineRule="auto"/><w:contextualSpacing/><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/><w:i/><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p>')
